$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddr, $text) {
    $rng = $ws.Range($rangeAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
}

Set-TextValue "D2" "625,910.56"
Set-TextValue "E2" "-1,974,089.44"
Set-TextValue "F2" "24.07%"
Set-TextValue "H2" "50.00"
Set-TextValue "I2" "-130.00"
Set-TextValue "J2" "27.78%"

Set-TextValue "D3" "232,450.83"
Set-TextValue "E3" "-967,549.17"
Set-TextValue "F3" "19.37%"
Set-TextValue "H3" "75.00"
Set-TextValue "I3" "-163.00"
Set-TextValue "J3" "31.51%"

Set-TextValue "D4" "319,917.67"
Set-TextValue "E4" "-1,680,082.33"
Set-TextValue "F4" "16.00%"
Set-TextValue "H4" "86.00"
Set-TextValue "I4" "-142.00"
Set-TextValue "J4" "37.72%"

Set-TextValue "D5" "235,262.17"
Set-TextValue "E5" "-1,064,737.83"
Set-TextValue "F5" "18.10%"
Set-TextValue "H5" "72.00"
Set-TextValue "I5" "-138.00"
Set-TextValue "J5" "34.29%"

Set-TextValue "D6" "253,032.96"
Set-TextValue "E6" "-1,146,967.04"
Set-TextValue "F6" "18.07%"
Set-TextValue "H6" "64.00"
Set-TextValue "I6" "-136.00"
Set-TextValue "J6" "32.00%"

Set-TextValue "D7" "1,666,574.19"
Set-TextValue "E7" "-6,833,425.81"
Set-TextValue "F7" "19.12%"
Set-TextValue "H7" "347.00"
Set-TextValue "I7" "-709.00"
Set-TextValue "J7" "32.66%"
